# Update the "取得日時" (retrieved datetime) timestamp in column A for the
# data rows on the "ランサーズ" sheet from 2025-12-16 12:39:47 to
# 2025-12-16 12:53:23 (appended a new batch of rows at 12:53 JST, but the
# underlying diff only shows the existing rows' timestamp text changing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-16 12:39:47"
$newTimestamp = "2025-12-16 12:53:23"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
